$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2518.2144
$ws.Range("J19").Value = 3700.5715
$ws.Range("L19").Value = 3700.5715
$ws.Range("N19").Value = -4050.5715
$ws.Range("H28").Value = 1753.7826
$ws.Range("I28").Value = 525.9286
$ws.Range("K28").Value = 525.9286
$ws.Range("M28").Value = -40.92859999999996
$ws.Range("H53").Value = 338.76923
$ws.Range("I53").Value = 237.71428
$ws.Range("J53").Value = 456.66666
$ws.Range("K53").Value = 237.71428
$ws.Range("L53").Value = 456.66666
$ws.Range("M53").Value = 399.28572
$ws.Range("N53").Value = -1730.66666
$ws.Range("H86").Value = 4486.8887
$ws.Range("I86").Value = 3745.5715
$ws.Range("K86").Value = 3745.5715
$ws.Range("M86").Value = -2622.5715
$ws.Range("H89").Value = 4486.8887
$ws.Range("I89").Value = 3745.5715
$ws.Range("K89").Value = 18727.8575
$ws.Range("M89").Value = -13111.8575
$ws.Range("H97").Value = 2208.2222
$ws.Range("J97").Value = 2208.2222
$ws.Range("L97").Value = 6624.6666
$ws.Range("N97").Value = -7616.6666
$ws.Range("H135").Value = 746.3913
$ws.Range("I135").Value = 464.64706
$ws.Range("K135").Value = 4181.82354
$ws.Range("M135").Value = -1646.82354

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 66663.336
$ws.Range("J92").Value = 66663.336
$ws.Range("L92").Value = 66663.336
$ws.Range("N92").Value = -71655.336
$ws.Range("H94").Value = 65000
$ws.Range("J94").Value = 65000
$ws.Range("L94").Value = 65000
$ws.Range("N94").Value = -66802
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H132").Value = 32305766
$ws.Range("I132").Value = 2547.2593
$ws.Range("J132").Value = 250352500
$ws.Range("K132").Value = 7641.777900000001
$ws.Range("L132").Value = 751057500
$ws.Range("M132").Value = -5111.777900000001
$ws.Range("N132").Value = -751062560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18802.059
$ws.Range("I86").Value = 7843.3335
$ws.Range("K86").Value = 7843.3335
$ws.Range("M86").Value = -6720.3335
$ws.Range("H89").Value = 18802.059
$ws.Range("I89").Value = 7843.3335
$ws.Range("K89").Value = 39216.6675
$ws.Range("M89").Value = -33600.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 659.8570999999999
$ws.Range("I16").Value = 659.8570999999999
$ws.Range("K16").Value = 659.8570999999999
$ws.Range("M16").Value = -372.8570999999999
$ws.Range("H88").Value = 38331.668
$ws.Range("J88").Value = 38331.668
$ws.Range("L88").Value = 38331.668
$ws.Range("N88").Value = -39143.668
$ws.Range("H91").Value = 38331.668
$ws.Range("J91").Value = 38331.668
$ws.Range("L91").Value = 38331.668
$ws.Range("N91").Value = -41139.668
$ws.Range("H99").Value = 6734.8423
$ws.Range("I99").Value = 6692.5
$ws.Range("K99").Value = 6692.5
$ws.Range("M99").Value = -5194.5
$ws.Range("H113").Value = 659.8570999999999
$ws.Range("I113").Value = 659.8570999999999
$ws.Range("K113").Value = 659.8570999999999
$ws.Range("M113").Value = 1510.1429
$ws.Range("H115").Value = 67165.5
$ws.Range("I115").Value = 64333
$ws.Range("J115").Value = 69998
$ws.Range("K115").Value = 64333
$ws.Range("L115").Value = 69998
$ws.Range("M115").Value = -63158
$ws.Range("N115").Value = -72348
$ws.Range("H126").Value = 6734.8423
$ws.Range("I126").Value = 6692.5
$ws.Range("K126").Value = 20077.5
$ws.Range("M126").Value = -17607.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42194700
$ws.Range("I4").Value = 57356892
$ws.Range("K4").Value = 172070676
$ws.Range("M4").Value = -172070564
$ws.Range("H129").Value = 4184
$ws.Range("J129").Value = 3832.95
$ws.Range("L129").Value = 11498.85
$ws.Range("N129").Value = -21498.85
$ws.Range("H132").Value = 1566.875
$ws.Range("I132").Value = 1548.25
$ws.Range("K132").Value = 13934.25
$ws.Range("M132").Value = -11404.25
$ws.Range("H136").Value = 2333
$ws.Range("I136").Value = 2333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1899
$ws.Range("N136").Value = $null
$ws.Range("H140").Value = 2732.875
$ws.Range("I140").Value = 2632.8333
$ws.Range("K140").Value = 7898.499899999999
$ws.Range("M140").Value = -2718.499899999999
$ws.Range("H141").Value = 7862.76
$ws.Range("J141").Value = 11156.5
$ws.Range("L141").Value = 33469.5
$ws.Range("N141").Value = -43829.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 30000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 30000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -30490
$ws.Range("H80").Value = 4864.222
$ws.Range("I80").Value = 4721.5
$ws.Range("J80").Value = 6006
$ws.Range("K80").Value = 4721.5
$ws.Range("L80").Value = 6006
$ws.Range("M80").Value = -3723.5
$ws.Range("N80").Value = -8002
$ws.Range("H83").Value = 4864.222
$ws.Range("I83").Value = 4721.5
$ws.Range("J83").Value = 6006
$ws.Range("K83").Value = 23607.5
$ws.Range("L83").Value = 30030
$ws.Range("M83").Value = -18615.5
$ws.Range("N83").Value = -40014
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H133").Value = 148000
$ws.Range("J133").Value = 148000
$ws.Range("L133").Value = 148000
$ws.Range("N133").Value = -158120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2804.4443
$ws.Range("I22").Value = 1397.2858
$ws.Range("K22").Value = 1397.2858
$ws.Range("M22").Value = -1102.2858
$ws.Range("H27").Value = 2804.4443
$ws.Range("I27").Value = 1397.2858
$ws.Range("K27").Value = 1397.2858
$ws.Range("M27").Value = -1290.2858
$ws.Range("H40").Value = 4168.3687
$ws.Range("I40").Value = 4312
$ws.Range("K40").Value = 4312
$ws.Range("M40").Value = -4176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3314.5
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 3721.75
$ws.Range("K81").Value = 5000
$ws.Range("L81").Value = 7443.5
$ws.Range("M81").Value = -3939
$ws.Range("N81").Value = -9565.5
$ws.Range("H84").Value = 3314.5
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 3721.75
$ws.Range("K84").Value = 25000
$ws.Range("L84").Value = 37217.5
$ws.Range("M84").Value = -19696
$ws.Range("N84").Value = -47825.5
$ws.Range("H113").Value = 716.7
$ws.Range("I113").Value = 381
$ws.Range("K113").Value = 1143
$ws.Range("M113").Value = 1027
$ws.Range("H122").Value = 4748.5
$ws.Range("I122").Value = 4497.5
$ws.Range("K122").Value = 13492.5
$ws.Range("M122").Value = -11042.5
